# Added matches from May 4th, 2023.
# The ELO ratings were recomputed after the new matches were played,
# which reshuffled several players' rank order in the table and updated
# their current_elo values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$playerNames = @(
    'Aaron Carter'
    'Will Simpson'
    'Levin Lee'
    'Rohan Chowla'
    'Kevin Cooper'
    'Roman Ramirez'
    'Gabe Silverstein'
    'Jack Massingill'
    'Nathan Snow'
    'Kristian Banlaoi'
    'Yvonne Nguyen'
    'Coby Lovelace'
    'Ann Hall'
    'Piper Parker'
    'Carla Betancourt'
    'Cason Duszak'
    'Reagan Fryatt'
    'Anna Brown'
    'Rose Roché'
    'Helen Dunn'
    'Noah Dale'
    'Jason Jackson'
    'Brian Tafazoli'
    'Matthew Rusten'
    'Leah Baetcke'
    'Evan Sooklal'
    'Sam Carswell-Tellis'
    'Cassie Deering'
    'Paul Bartenfeld'
)

$eloValues = @(
    1537.121990870036
    1396.331533632466
    1391.624867430988
    1277.890315082629
    1271.813139769939
    1270.649156167338
    1233.467836182129
    1228.341892966386
    1219.740685506919
    1216.037581761179
    1204.868130393788
    1198.883494704046
    1198.228992509007
    1194.748819295017
    1192.258245275334
    1183.651911226616
    1182.605339128595
    1176.279801661901
    1145.11049623546
    1143.988590686295
    1143.988590686295
    1140.220115271984
    1140
    1139.327236294392
    1138.692457411543
    1113.066776939508
    1076.165195726573
    1058.266119428515
    986.6306877551247
)

for ($i = 0; $i -lt $playerNames.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $playerNames[$i]
    $ws.Cells.Item($row, 2).Value = $eloValues[$i]
}

